$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: "Size (GB)" ---
# Copy the header formatting from an existing header cell (bold, centered,
# bordered) onto E1, then set its text.
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E1").Value = "Size (GB)"

# Blank placeholder cells for the existing rows (column E had no data for
# these datasets).
$ws.Range("E2").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("E7").Value = ""

# --- New rows for the Gasperini/Shendure datasets ---
$ws.Range("A8").Value = "GasperiniShendure2019_atscale.h5ad"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = 1.738121328875422

$ws.Range("A9").Value = "GasperiniShendure2019_highMOI.h5ad"
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = 0.3896072432398796

$ws.Range("A10").Value = "GasperiniShendure2019_lowMOI.h5ad"
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = 0.300553466193378
